# Apply updated capacity values on the "Capacità di trasmissione MW" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Capacità di trasmissione MW")

$ws.Range("C2").Value = 5500

$ws.Range("B3").Value = 4700
$ws.Range("D3").Value = 4600
$ws.Range("H3").Value = 400

$ws.Range("C4").Value = 4500
$ws.Range("E4").Value = 2900
$ws.Range("G4").Value = 1000

$ws.Range("D5").Value = 5700

$ws.Range("E6").Value = 3000
$ws.Range("G6").Value = 2100

$ws.Range("D7").Value = 1000
$ws.Range("F7").Value = 2000
$ws.Range("H7").Value = 1000

$ws.Range("C8").Value = 400
$ws.Range("G8").Value = 1000

$ws.Range("F9").Value = 1700
